# Daily attendance processing - 2025-12-21 10:28:57
# Reorders the comma-separated "Recorded By" (column G) values for specific rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 'system, System, backup@backdoor.com'
$ws.Cells.Item(3, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(4, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(6, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(10, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(11, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(12, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(13, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(14, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(15, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(17, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(18, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(19, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(20, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(21, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(22, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(24, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(26, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(28, 7).Value = 'system, System, backup@backdoor.com'
$ws.Cells.Item(29, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(30, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(32, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(36, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(37, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(38, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(39, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(40, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(41, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(43, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(44, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(45, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(46, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(47, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(48, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(50, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(52, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(54, 7).Value = 'system, System, backup@backdoor.com'
$ws.Cells.Item(55, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(56, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(58, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(62, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(63, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(64, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(65, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(66, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(67, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(69, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(70, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(71, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(72, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(73, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(74, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(76, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(78, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(83, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(84, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(85, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(86, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(90, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(92, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(93, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(94, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(96, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(99, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(101, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(109, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(110, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(111, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(112, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(116, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(118, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(119, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(120, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(122, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(125, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(127, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(135, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(136, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(137, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(138, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(142, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(144, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(145, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(146, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(148, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(151, 7).Value = 'dnasr281@gmail.com, System'
$ws.Cells.Item(153, 7).Value = 'dnasr281@gmail.com, System'
